$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1 - copy formatting (style) from H1, then set the text values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I (I0) and J (IF) columns, rows 2-20
$data = @(
    @(2,  6,  7),
    @(3,  5,  6),
    @(4,  2,  3),
    @(5,  6,  6),
    @(6,  8,  9),
    @(7,  3,  5),
    @(8,  9,  9),
    @(9,  4,  4),
    @(10, 6,  6),
    @(11, 8,  9),
    @(12, 10, 11),
    @(13, 2,  3),
    @(14, 8,  8),
    @(15, 8,  8),
    @(16, 5,  5),
    @(17, 8,  8),
    @(18, 3,  3),
    @(19, 6,  6),
    @(20, 6,  6)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
